# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gains a new (empty) column inserted
# before the existing "Late" column (old column N), pushing the old
# N/O/P ("Late", "Outstanding heading", "Outstanding") out to O/P/Q.
# The "Repayment schedule" sheet becomes the active sheet/tab (it was
# "Transactions" before), and the selection on "Repayment schedule"
# moves from H8 to the new corresponding cell S8.

$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (shifts old N/O/P -> O/P/Q).
$wsRepay.Columns("N:N").Insert()

# Give the newly inserted column its own (non bestFit) width, matching
# the width of the "In Advance" column (M) immediately to its left.
$wsRepay.Columns("N:N").ColumnWidth = 9.8

# "Repayment schedule" becomes the active sheet/tab, with the
# selection moved to S8 (the cell that used to be H8 before the
# column insert shifted everything right of M over by one).
$wsRepay.Activate()
$wsRepay.Range("S8").Select()
